# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    # Force the literal text into the cell without Excel's automatic
    # number/date coercion (e.g. "449.80" -> 449.8), then drop the
    # temporary text format so the cell's style stays untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "62.408.15"
Set-TextValue $ws.Range("E2") "  +0.12%  "
Set-TextValue $ws.Range("D3") "3.013.28"
Set-TextValue $ws.Range("E3") "  +0.17%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "594.08"
Set-TextValue $ws.Range("E5") "  +2.02%  "
Set-TextValue $ws.Range("D6") "147.59"
Set-TextValue $ws.Range("E6") "  +1.29%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("D8") "3.009.57"
Set-TextValue $ws.Range("E8") "  +0.08%  "
Set-TextValue $ws.Range("D9") "0.517"
Set-TextValue $ws.Range("E9") "  -2.06%  "
Set-TextValue $ws.Range("E10") "  +8.23%  "
Set-TextValue $ws.Range("D11") "0.149"
Set-TextValue $ws.Range("E11") "  +0.52%  "
Set-TextValue $ws.Range("E12") "  -1.65%  "
Set-TextValue $ws.Range("D13") "0.0000232"
Set-TextValue $ws.Range("E13") "  +1.51%  "
Set-TextValue $ws.Range("D14") "34.43"
Set-TextValue $ws.Range("E14") "  -0.11%  "
Set-TextValue $ws.Range("E15") "  +2.40%  "
Set-TextValue $ws.Range("D16") "3.513.55"
Set-TextValue $ws.Range("E16") "  +0.38%  "
Set-TextValue $ws.Range("D17") "62.264.24"
Set-TextValue $ws.Range("E17") "  -0.03%  "
Set-TextValue $ws.Range("D18") "6.99"
Set-TextValue $ws.Range("E18") "  -1.76%  "
Set-TextValue $ws.Range("D19") "3.019.56"
Set-TextValue $ws.Range("E19") "  +0.44%  "
Set-TextValue $ws.Range("D20") "449.80"
Set-TextValue $ws.Range("E20") "  -1.34%  "
Set-TextValue $ws.Range("E21") "  +1.16%  "
Set-TextValue $ws.Range("D22") "0.686"
Set-TextValue $ws.Range("E22") "  -0.32%  "
Set-TextValue $ws.Range("D23") "7.36"
Set-TextValue $ws.Range("E23") "  -0.51%  "
Set-TextValue $ws.Range("D24") "81.96"
Set-TextValue $ws.Range("E24") "  +0.29%  "
Set-TextValue $ws.Range("D25") "11.07"
Set-TextValue $ws.Range("E25") "  +10.64%  "
Set-TextValue $ws.Range("D26") "2.25"
Set-TextValue $ws.Range("E26") "  +1.61%  "
Set-TextValue $ws.Range("D27") "12.17"
Set-TextValue $ws.Range("E27") "  -1.53%  "
Set-TextValue $ws.Range("E28") "  -0.01%  "
Set-TextValue $ws.Range("D29") "2.70"
Set-TextValue $ws.Range("E29") "  +3.57%  "
Set-TextValue $ws.Range("D30") "7.30"
Set-TextValue $ws.Range("E30") "  +4.80%  "
Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  +0.02%  "
Set-TextValue $ws.Range("D32") "2.09"
Set-TextValue $ws.Range("E32") "  +0.30%  "
Set-TextValue $ws.Range("D33") "27.40"
Set-TextValue $ws.Range("E34") "  +1.54%  "
Set-TextValue $ws.Range("D35") "0.0₃0849"
Set-TextValue $ws.Range("E35") "  +6.24%  "
Set-TextValue $ws.Range("E36") "  -0.29%  "
Set-TextValue $ws.Range("D37") "5.82"
Set-TextValue $ws.Range("E37") "  +0.98%  "
Set-TextValue $ws.Range("D38") "50.30"
Set-TextValue $ws.Range("E38") "  +0.11%  "
Set-TextValue $ws.Range("D39") "2.06"
Set-TextValue $ws.Range("E39") "  -2.73%  "
Set-TextValue $ws.Range("D40") "2.97"
Set-TextValue $ws.Range("E40") "  +2.81%  "
Set-TextValue $ws.Range("D41") "8.96"
Set-TextValue $ws.Range("E41") "  -2.26%  "
Set-TextValue $ws.Range("E42") "  +6.76%  "
Set-TextValue $ws.Range("D43") "402.82"
Set-TextValue $ws.Range("E43") "  +3.07%  "
Set-TextValue $ws.Range("D44") "41.10"
Set-TextValue $ws.Range("E44") "  +10.79%  "
Set-TextValue $ws.Range("D45") "0.278"
Set-TextValue $ws.Range("E45") "  +3.67%  "
Set-TextValue $ws.Range("D46") "0.0353"
Set-TextValue $ws.Range("E46") "  -1.25%  "
Set-TextValue $ws.Range("D47") "2.717.20"
Set-TextValue $ws.Range("E47") "  -0.21%  "
Set-TextValue $ws.Range("D48") "132.80"
Set-TextValue $ws.Range("E48") "  +3.32%  "
Set-TextValue $ws.Range("E49") "  +0.09%  "
Set-TextValue $ws.Range("D50") "2.19"
Set-TextValue $ws.Range("E50") "  -0.66%  "
Set-TextValue $ws.Range("D51") "0.107"
Set-TextValue $ws.Range("E51") "  -1.61%  "
